$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test credentials (Usuario / NroSiniestro) used for the inspection
$ws.Range("D2").Value = "tcorvetto"
$ws.Range("F2").Value = "'0420172008282"

# Move active selection to F3 as in the updated workbook
$ws.Range("F3").Select()
